# Auto-generated PowerShell/COM script implementing the "fix SF3B1alpha conditions"
# commit: adds 3 new example rows (zeroScore, withNAs, SF3B1alphaTestCase) plus
# trailing blank-but-styled helper rows/columns used by the sheet formulas, and
# the 5 new shared strings those rows introduce.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet  # "examples" sheet is the active/selected one in this workbook

# --- 0) Seed the 5 brand-new shared strings in the same first-use order as the
#        authored workbook (withNAs, SF3B1alphaTestCase, Intermediate, zeroScore,
#        "Very Good") so the shared-strings table lines up index-for-index. ---
$ws.Range("A5").Value2 = "withNAs"
$ws.Range("A6").Value2 = "SF3B1alphaTestCase"
$ws.Range("J6").Value2 = "Intermediate"
$ws.Range("A4").Value2 = "zeroScore"
$ws.Range("J4").Value2 = "Very Good"

# --- 1) Write the literal values for the 3 new data rows (4,5,6) ---
# Row 4
$ws.Range("A4").Value2 = "zeroScore"
$ws.Range("B4").Value2 = 10
$ws.Range("C4").Value2 = 10
$ws.Range("D4").Value2 = 10
$ws.Range("E4").Value2 = 0
$ws.Range("F4").Value2 = 100
$ws.Range("G4").Value2 = 0
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = "Very Good"
$ws.Range("K4").Value2 = 0
$ws.Range("L4").Value2 = 0
$ws.Range("M4").Value2 = 0
$ws.Range("N4").Value2 = 0
$ws.Range("O4").Value2 = 0
$ws.Range("P4").Value2 = 0
$ws.Range("Q4").Value2 = 0
$ws.Range("R4").Value2 = 1
$ws.Range("S4").Value2 = 0
$ws.Range("T4").Value2 = 0
$ws.Range("U4").Value2 = 0
$ws.Range("V4").Value2 = 0
$ws.Range("W4").Value2 = 0
$ws.Range("X4").Value2 = 0
$ws.Range("Y4").Value2 = 0
$ws.Range("Z4").Value2 = 0
$ws.Range("AA4").Value2 = 0
$ws.Range("AB4").Value2 = 0
$ws.Range("AC4").Value2 = 0
$ws.Range("AD4").Value2 = 0
$ws.Range("AE4").Value2 = 0
$ws.Range("AF4").Value2 = 0
$ws.Range("AG4").Value2 = 0
$ws.Range("AH4").Value2 = 0
$ws.Range("AI4").Value2 = 0
$ws.Range("AJ4").Value2 = 0
$ws.Range("AK4").Value2 = 0
$ws.Range("AL4").Value2 = 0
$ws.Range("AM4").Value2 = 0
$ws.Range("AN4").Value2 = 0
$ws.Range("AO4").Value2 = 0
$ws.Range("AP4").Value2 = 0
$ws.Range("AQ4").Value2 = 0
$ws.Range("AR4").Value2 = 0

# Row 5
$ws.Range("A5").Value2 = "withNAs"
$ws.Range("B5").Value2 = 9.6
$ws.Range("C5").Value2 = 281
$ws.Range("D5").Value2 = 9
$ws.Range("E5").Value2 = 4.84
$ws.Range("F5").Value2 = 79
$ws.Range("G5").Value2 = 0
$ws.Range("H5").Value2 = 0
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = "Good"
$ws.Range("K5").Value2 = 0
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = "NA"
$ws.Range("N5").Value2 = 0
$ws.Range("O5").Value2 = 0
$ws.Range("P5").Value2 = 0
$ws.Range("Q5").Value2 = 1
$ws.Range("R5").Value2 = "NA"
$ws.Range("S5").Value2 = 0
$ws.Range("T5").Value2 = 0
$ws.Range("U5").Value2 = 0
$ws.Range("V5").Value2 = 0
$ws.Range("W5").Value2 = 0
$ws.Range("X5").Value2 = 0
$ws.Range("Y5").Value2 = 0
$ws.Range("Z5").Value2 = 0
$ws.Range("AA5").Value2 = 0
$ws.Range("AB5").Value2 = 0
$ws.Range("AC5").Value2 = 0
$ws.Range("AD5").Value2 = 0
$ws.Range("AE5").Value2 = 0
$ws.Range("AF5").Value2 = 0
$ws.Range("AG5").Value2 = 0
$ws.Range("AH5").Value2 = 1
$ws.Range("AI5").Value2 = 1
$ws.Range("AJ5").Value2 = 0
$ws.Range("AK5").Value2 = "NA"
$ws.Range("AL5").Value2 = 0
$ws.Range("AM5").Value2 = "NA"
$ws.Range("AN5").Value2 = 0
$ws.Range("AO5").Value2 = 0
$ws.Range("AP5").Value2 = "NA"
$ws.Range("AQ5").Value2 = 0
$ws.Range("AR5").Value2 = "NA"

# Row 6
$ws.Range("A6").Value2 = "SF3B1alphaTestCase"
$ws.Range("B6").Value2 = 10
$ws.Range("C6").Value2 = 100
$ws.Range("D6").Value2 = 5
$ws.Range("E6").Value2 = "NA"
$ws.Range("F6").Value2 = "NA"
$ws.Range("G6").Value2 = 0
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 0
$ws.Range("J6").Value2 = "Intermediate"
$ws.Range("K6").Value2 = 0
$ws.Range("L6").Value2 = 0
$ws.Range("M6").Value2 = "NA"
$ws.Range("N6").Value2 = 0
$ws.Range("O6").Value2 = 0
$ws.Range("P6").Value2 = 0
$ws.Range("Q6").Value2 = 0
$ws.Range("R6").Value2 = 0
$ws.Range("S6").Value2 = 0
$ws.Range("T6").Value2 = 0
$ws.Range("U6").Value2 = 0
$ws.Range("V6").Value2 = 0
$ws.Range("W6").Value2 = 0
$ws.Range("X6").Value2 = 0
$ws.Range("Y6").Value2 = 0
$ws.Range("Z6").Value2 = 0
$ws.Range("AA6").Value2 = 0
$ws.Range("AB6").Value2 = 0
$ws.Range("AC6").Value2 = 0
$ws.Range("AD6").Value2 = 0
$ws.Range("AE6").Value2 = "NA"
$ws.Range("AF6").Value2 = 0
$ws.Range("AG6").Value2 = 0
$ws.Range("AH6").Value2 = 0
$ws.Range("AI6").Value2 = 0
$ws.Range("AJ6").Value2 = 0
$ws.Range("AK6").Value2 = 0
$ws.Range("AL6").Value2 = 0
$ws.Range("AM6").Value2 = 0
$ws.Range("AN6").Value2 = 0
$ws.Range("AO6").Value2 = 0
$ws.Range("AP6").Value2 = 0
$ws.Range("AQ6").Value2 = 0
$ws.Range("AR6").Value2 = 0

# --- 2) Re-apply the "highlight" style (s=9, same as column A/B/J/M/AK.. in the
#        existing rows) to the fixed set of columns that carry it on every data
#        row, plus stamp that same style (with no value) on rows 7-13 so the
#        sheet's helper grid extends all the way down, matching the template.
$styleCols = @("A","B","J","M","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","BB","BE","CC","CD","CE","CF","CG","CH","CI","CJ","CK","CL","CT","CW")
$ws.Range("A2").Copy() | Out-Null
for ($r = 4; $r -le 13; $r++) {
    foreach ($col in $styleCols) {
        $addr = $col + $r
        $ws.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
    }
}
$excel.CutCopyMode = 0

# --- 3) Restore the cursor/selection like the saved workbook (cell V12) ---
$ws.Range("V12").Select() | Out-Null

